# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (number of people interested) in column F
# for worksheets "展览" (sheet 1) and "全部类型" (sheet 4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 250
$ws1.Range("F5").Value  = 2968
$ws1.Range("F6").Value  = 2006
$ws1.Range("F7").Value  = 387
$ws1.Range("F8").Value  = 135
$ws1.Range("F9").Value  = 1100
$ws1.Range("F10").Value = 204
$ws1.Range("F11").Value = 593
$ws1.Range("F12").Value = 60

# --- Sheet "全部类型" (index 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 250
$ws4.Range("F5").Value  = 2968
$ws4.Range("F6").Value  = 2006
$ws4.Range("F7").Value  = 387
$ws4.Range("F9").Value  = 135
$ws4.Range("F10").Value = 1100
$ws4.Range("F11").Value = 204
$ws4.Range("F12").Value = 593
$ws4.Range("F13").Value = 60
